# Sara-Alert-Format-Invalid-Fields.xlsx fixture update:
# - zoom the sheet view to 91%
# - change the "Preferred Contact Time" (AW) column test values on rows 3-7
# - add two new sample rows (6 & 7) mirroring row 5's monitoree data
# - drop the leftover blank styled cells from row 7

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- sheet view zoom ---
$excel.ActiveWindow.Zoom = 91

# --- Rows 6 & 7: mirror row 5's monitoree sample data (copy values only,
# so text like "false"/"true" stays text instead of becoming a boolean),
# then drop the handful of columns row 5 has that rows 6/7 should not.
$ws.Range("A5:AT5").Copy()
$ws.Range("A6:AT6").PasteSpecial(-4163)
$ws.Range("A5:AT5").Copy()
$ws.Range("A7:AT7").PasteSpecial(-4163)

foreach ($r in 6, 7) {
    foreach ($col in "B", "D", "U", "AF", "AG", "AH", "AJ", "AR") {
        $ws.Range($col + $r).Clear()
    }
}

# --- AW column edge-case values (numbers with new numeric formats) ---
# AW5: was shared-string "24:00" -> becomes literal 0 formatted "0.00"
$ws.Range("AW5").NumberFormat = "0.00"
$ws.Range("AW5").Value = 0

# AW4: was a big literal number -> becomes 16 formatted "0.00E+00"
$ws.Range("AW4").NumberFormat = "0.00E+00"
$ws.Range("AW4").Value = 16

# AW6 (new row): literal 1 formatted "[h]:mm:ss"
$ws.Range("AW6").NumberFormat = "[h]:mm:ss"
$ws.Range("AW6").Value = 1

# AW7 (new row): new shared string "15"
$ws.Range("AW7").Value = "15"

# AW3 keeps referencing the "1:30" text (shared string table gets
# compacted once "24:00" stops being referenced above)
$ws.Range("AW3").Value = "1:30"

# --- drop the stray blank styled cells that used to live on row 7 ---
$ws.Range("D7").Clear()
$ws.Range("BL7").Clear()
$ws.Range("BM7").Clear()
$ws.Range("BO7").Clear()
$ws.Range("CH7").Clear()
$ws.Range("CO7").Clear()
$ws.Range("CP7").Clear()
